# The source document carries a handful of legacy SharePoint "Document
# Information Panel" custom XML parts (the content-type schema, the
# SharePoint form-template pointer, and the empty documentManagement
# properties bag, plus their accompanying itemProps datastore items).
# They aren't bound to anything in the document body (no XML mapped
# content controls reference them) - they are just inert metadata left
# over from the last time the file was checked in/out of a SharePoint
# document library. This "Update the Documentserver Docker details"
# housekeeping pass drops that dead metadata so the package only keeps
# the parts the document actually uses.
#
# Remove every custom XML part whose namespace matches the SharePoint /
# Document Information Panel metadata schemas. This mirrors what
# Document.CustomXMLParts / CustomXMLPart.Delete do in real Word
# automation: enumerate the parts carrying this metadata and delete
# them one by one.

$d = $word.ActiveDocument

$sharePointNamespaces = @(
    "http://schemas.microsoft.com/office/2006/metadata/contentType",
    "http://schemas.microsoft.com/sharepoint/v3/contenttype/forms",
    "http://schemas.microsoft.com/office/2006/metadata/properties"
)

ForEach ($ns in $sharePointNamespaces) {
    $scoped = $d.CustomXMLParts.SelectByNamespace($ns)
    For ($i = $scoped.Count; $i -ge 1; $i--) {
        $scoped.Item($i).Delete()
    }
}

# Safety net: walk whatever is left back-to-front and delete any
# remaining non-built-in part (covers hosts where SelectByNamespace
# doesn't resolve but the flat collection still does).
For ($i = $d.CustomXMLParts.Count; $i -ge 1; $i--) {
    $part = $d.CustomXMLParts.Item($i)
    If (-not $part.BuiltIn) {
        $part.Delete()
    }
}

Write-Output ("CustomXMLParts remaining: " + $d.CustomXMLParts.Count)
